# Update ranking at 2025-12-04 12:18
# Append a new tracking row (row 52) to the ranking sheet with the
# latest timestamp and placeholder "-" values for the ranking columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A52").Value = "2025/12/04 21:00"
$ws.Range("B52").Value = "-"
$ws.Range("C52").Value = "-"
$ws.Range("D52").Value = "-"
$ws.Range("E52").Value = "-"
$ws.Range("F52").Value = "-"
$ws.Range("G52").Value = "-"
